$wb = $excel.ActiveWorkbook

# Column G width: widen col 7 from 6.998 to 8.141 (engine quantizes to 1/6-units;
# 7.3333333333333335 is the ColumnWidth input that rounds closest to the target).
$targetColWidth = 7.3333333333333335

# --- Sheet "Sum total" ---
$ws = $wb.Worksheets.Item("Sum total")
$ws.Columns.Item(7).ColumnWidth = $targetColWidth
$ws.Range("G37").Value = 2200.0
$ws.Range("AI37").Value = -2200.0
$ws.Range("G38").Value = 2420.0
$ws.Range("AI38").Value = -2420.0
$ws.Range("G39").Value = 2662.0
$ws.Range("AI39").Value = -2662.0
$ws.Range("AJ39").Value = -2662.0
$ws.Range("G40").Value = 2928.0
$ws.Range("AI40").Value = -2928.0
$ws.Range("AJ40").Value = -5590.0
$ws.Range("G41").Value = 3221.0
$ws.Range("AI41").Value = -3221.0
$ws.Range("AJ41").Value = -8811.0
$ws.Range("G42").Value = 3543.0
$ws.Range("AI42").Value = -3543.0
$ws.Range("AJ42").Value = -12354.0
$ws.Range("G43").Value = 3898.0
$ws.Range("AI43").Value = -3898.0
$ws.Range("AJ43").Value = -16252.0
$ws.Range("G44").Value = 4287.0
$ws.Range("AI44").Value = -4287.0
$ws.Range("AJ44").Value = -20539.0
$ws.Range("G45").Value = 4716.0
$ws.Range("AI45").Value = -4716.0
$ws.Range("AJ45").Value = -25255.0
$ws.Range("G46").Value = 5188.0
$ws.Range("AI46").Value = -5188.0
$ws.Range("AJ46").Value = -30443.0
$ws.Range("G47").Value = 5706.0
$ws.Range("AI47").Value = -5706.0
$ws.Range("AJ47").Value = -36149.0
$ws.Range("G48").Value = 6277.0
$ws.Range("AI48").Value = -6277.0
$ws.Range("AJ48").Value = -42426.0
$ws.Range("G49").Value = 6904.0
$ws.Range("AI49").Value = -6904.0
$ws.Range("AJ49").Value = -49330.0
$ws.Range("G50").Value = 7595.0
$ws.Range("AI50").Value = -7595.0
$ws.Range("AJ50").Value = -56925.0
$ws.Range("G51").Value = 8355.0
$ws.Range("AI51").Value = -8355.0
$ws.Range("AJ51").Value = -65280.0
$ws.Range("G52").Value = 9190.0
$ws.Range("AI52").Value = -9190.0
$ws.Range("AJ52").Value = -74470.0
$ws.Range("G53").Value = 10109.0
$ws.Range("AI53").Value = -10109.0
$ws.Range("AJ53").Value = -84579.0
$ws.Range("G54").Value = 11120.0
$ws.Range("AI54").Value = -11120.0
$ws.Range("AJ54").Value = -95699.0
$ws.Range("G55").Value = 12232.0
$ws.Range("AI55").Value = -12232.0
$ws.Range("AJ55").Value = -107931.0
$ws.Range("G56").Value = 13455.0
$ws.Range("AI56").Value = -13455.0
$ws.Range("AJ56").Value = -121386.0
$ws.Range("G57").Value = 14801.0
$ws.Range("AI57").Value = -14801.0
$ws.Range("AJ57").Value = -136187.0
$ws.Range("G58").Value = 16281.0
$ws.Range("AI58").Value = -16281.0
$ws.Range("AJ58").Value = -152468.0
$ws.Range("G59").Value = 17909.0
$ws.Range("AI59").Value = -17909.0
$ws.Range("AJ59").Value = -170377.0
$ws.Range("G60").Value = 19699.0
$ws.Range("AI60").Value = -19699.0
$ws.Range("AJ60").Value = -190076.0
$ws.Range("G61").Value = 21670.0
$ws.Range("AI61").Value = -21670.0
$ws.Range("AJ61").Value = -211746.0
$ws.Range("G62").Value = 23836.0
$ws.Range("AI62").Value = -23836.0
$ws.Range("AJ62").Value = -235582.0
$ws.Range("G63").Value = 26220.0
$ws.Range("AI63").Value = -26220.0
$ws.Range("AJ63").Value = -261802.0
$ws.Range("G64").Value = 28842.0
$ws.Range("AI64").Value = -28842.0
$ws.Range("AJ64").Value = -290644.0
$ws.Range("G65").Value = 31726.0
$ws.Range("AI65").Value = -31726.0
$ws.Range("AJ65").Value = -322370.0
$ws.Range("G66").Value = 34899.0
$ws.Range("AI66").Value = -35348.44
$ws.Range("AJ66").Value = -357718.44
$ws.Range("G67").Value = 38389.0
$ws.Range("AI67").Value = -40583.38
$ws.Range("AJ67").Value = -398301.82
$ws.Range("G68").Value = 42228.0
$ws.Range("AI68").Value = -46341.81
$ws.Range("AJ68").Value = -444643.63
$ws.Range("G69").Value = 46450.0
$ws.Range("AI69").Value = -52675.2
$ws.Range("AJ69").Value = -497318.83
$ws.Range("G70").Value = 51095.0
$ws.Range("AI70").Value = -59642.72
$ws.Range("AJ70").Value = -556961.55
$ws.Range("G71").Value = 56205.0
$ws.Range("AI71").Value = -67307.49
$ws.Range("AJ71").Value = -624269.04
$ws.Range("G72").Value = 61826.0
$ws.Range("AI72").Value = -75738.73
$ws.Range("AJ72").Value = -700007.77

# --- Sheet "Sum private" ---
$ws = $wb.Worksheets.Item("Sum private")
$ws.Columns.Item(7).ColumnWidth = $targetColWidth
$ws.Range("G37").Value = 2200.0
$ws.Range("AI37").Value = -2200.0
$ws.Range("G38").Value = 2420.0
$ws.Range("AI38").Value = -2420.0
$ws.Range("G39").Value = 2662.0
$ws.Range("AI39").Value = -2662.0
$ws.Range("AJ39").Value = -2662.0
$ws.Range("G40").Value = 2928.0
$ws.Range("AI40").Value = -2928.0
$ws.Range("AJ40").Value = -5590.0
$ws.Range("G41").Value = 3221.0
$ws.Range("AI41").Value = -3221.0
$ws.Range("AJ41").Value = -8811.0
$ws.Range("G42").Value = 3543.0
$ws.Range("AI42").Value = -3543.0
$ws.Range("AJ42").Value = -12354.0
$ws.Range("G43").Value = 3898.0
$ws.Range("AI43").Value = -3898.0
$ws.Range("AJ43").Value = -16252.0
$ws.Range("G44").Value = 4287.0
$ws.Range("AI44").Value = -4287.0
$ws.Range("AJ44").Value = -20539.0
$ws.Range("G45").Value = 4716.0
$ws.Range("AI45").Value = -4716.0
$ws.Range("AJ45").Value = -25255.0
$ws.Range("G46").Value = 5188.0
$ws.Range("AI46").Value = -5188.0
$ws.Range("AJ46").Value = -30443.0
$ws.Range("G47").Value = 5706.0
$ws.Range("AI47").Value = -5706.0
$ws.Range("AJ47").Value = -36149.0
$ws.Range("G48").Value = 6277.0
$ws.Range("AI48").Value = -6277.0
$ws.Range("AJ48").Value = -42426.0
$ws.Range("G49").Value = 6904.0
$ws.Range("AI49").Value = -6904.0
$ws.Range("AJ49").Value = -49330.0
$ws.Range("G50").Value = 7595.0
$ws.Range("AI50").Value = -7595.0
$ws.Range("AJ50").Value = -56925.0
$ws.Range("G51").Value = 8355.0
$ws.Range("AI51").Value = -8355.0
$ws.Range("AJ51").Value = -65280.0
$ws.Range("G52").Value = 9190.0
$ws.Range("AI52").Value = -9190.0
$ws.Range("AJ52").Value = -74470.0
$ws.Range("G53").Value = 10109.0
$ws.Range("AI53").Value = -10109.0
$ws.Range("AJ53").Value = -84579.0
$ws.Range("G54").Value = 11120.0
$ws.Range("AI54").Value = -11120.0
$ws.Range("AJ54").Value = -95699.0
$ws.Range("G55").Value = 12232.0
$ws.Range("AI55").Value = -12232.0
$ws.Range("AJ55").Value = -107931.0
$ws.Range("G56").Value = 13455.0
$ws.Range("AI56").Value = -13455.0
$ws.Range("AJ56").Value = -121386.0
$ws.Range("G57").Value = 14801.0
$ws.Range("AI57").Value = -14801.0
$ws.Range("AJ57").Value = -136187.0
$ws.Range("G58").Value = 16281.0
$ws.Range("AI58").Value = -16281.0
$ws.Range("AJ58").Value = -152468.0
$ws.Range("G59").Value = 17909.0
$ws.Range("AI59").Value = -17909.0
$ws.Range("AJ59").Value = -170377.0
$ws.Range("G60").Value = 19699.0
$ws.Range("AI60").Value = -19699.0
$ws.Range("AJ60").Value = -190076.0
$ws.Range("G61").Value = 21670.0
$ws.Range("AI61").Value = -21670.0
$ws.Range("AJ61").Value = -211746.0
$ws.Range("G62").Value = 23836.0
$ws.Range("AI62").Value = -23836.0
$ws.Range("AJ62").Value = -235582.0
$ws.Range("G63").Value = 26220.0
$ws.Range("AI63").Value = -26220.0
$ws.Range("AJ63").Value = -261802.0
$ws.Range("G64").Value = 28842.0
$ws.Range("AI64").Value = -28842.0
$ws.Range("AJ64").Value = -290644.0
$ws.Range("G65").Value = 31726.0
$ws.Range("AI65").Value = -31726.0
$ws.Range("AJ65").Value = -322370.0
$ws.Range("G66").Value = 34899.0
$ws.Range("AI66").Value = -35348.44
$ws.Range("AJ66").Value = -357718.44
$ws.Range("G67").Value = 38389.0
$ws.Range("AI67").Value = -40583.38
$ws.Range("AJ67").Value = -398301.82
$ws.Range("G68").Value = 42228.0
$ws.Range("AI68").Value = -46341.81
$ws.Range("AJ68").Value = -444643.63
$ws.Range("G69").Value = 46450.0
$ws.Range("AI69").Value = -52675.2
$ws.Range("AJ69").Value = -497318.83
$ws.Range("G70").Value = 51095.0
$ws.Range("AI70").Value = -59642.72
$ws.Range("AJ70").Value = -556961.55
$ws.Range("G71").Value = 56205.0
$ws.Range("AI71").Value = -67307.49
$ws.Range("AJ71").Value = -624269.04
$ws.Range("G72").Value = 61826.0
$ws.Range("AI72").Value = -75738.73
$ws.Range("AJ72").Value = -700007.77

# --- Sheet "Cash" (also gets 22% tax-rate in column H) ---
$ws = $wb.Worksheets.Item("Cash")
$ws.Columns.Item(7).ColumnWidth = $targetColWidth
$ws.Range("G37").Value = 2200.0
$ws.Range("H37").Value = 0.22
$ws.Range("AI37").Value = -2200.0
$ws.Range("G38").Value = 2420.0
$ws.Range("H38").Value = 0.22
$ws.Range("AI38").Value = -2420.0
$ws.Range("G39").Value = 2662.0
$ws.Range("H39").Value = 0.22
$ws.Range("AI39").Value = -2662.0
$ws.Range("AJ39").Value = -2662.0
$ws.Range("G40").Value = 2928.0
$ws.Range("H40").Value = 0.22
$ws.Range("AI40").Value = -2928.0
$ws.Range("AJ40").Value = -5590.0
$ws.Range("G41").Value = 3221.0
$ws.Range("H41").Value = 0.22
$ws.Range("AI41").Value = -3221.0
$ws.Range("AJ41").Value = -8811.0
$ws.Range("G42").Value = 3543.0
$ws.Range("H42").Value = 0.22
$ws.Range("AI42").Value = -3543.0
$ws.Range("AJ42").Value = -12354.0
$ws.Range("G43").Value = 3898.0
$ws.Range("H43").Value = 0.22
$ws.Range("AI43").Value = -3898.0
$ws.Range("AJ43").Value = -16252.0
$ws.Range("G44").Value = 4287.0
$ws.Range("H44").Value = 0.22
$ws.Range("AI44").Value = -4287.0
$ws.Range("AJ44").Value = -20539.0
$ws.Range("G45").Value = 4716.0
$ws.Range("H45").Value = 0.22
$ws.Range("AI45").Value = -4716.0
$ws.Range("AJ45").Value = -25255.0
$ws.Range("G46").Value = 5188.0
$ws.Range("H46").Value = 0.22
$ws.Range("AI46").Value = -5188.0
$ws.Range("AJ46").Value = -30443.0
$ws.Range("G47").Value = 5706.0
$ws.Range("H47").Value = 0.22
$ws.Range("AI47").Value = -5706.0
$ws.Range("AJ47").Value = -36149.0
$ws.Range("G48").Value = 6277.0
$ws.Range("H48").Value = 0.22
$ws.Range("AI48").Value = -6277.0
$ws.Range("AJ48").Value = -42426.0
$ws.Range("G49").Value = 6904.0
$ws.Range("H49").Value = 0.22
$ws.Range("AI49").Value = -6904.0
$ws.Range("AJ49").Value = -49330.0
$ws.Range("G50").Value = 7595.0
$ws.Range("H50").Value = 0.22
$ws.Range("AI50").Value = -7595.0
$ws.Range("AJ50").Value = -56925.0
$ws.Range("G51").Value = 8355.0
$ws.Range("H51").Value = 0.22
$ws.Range("AI51").Value = -8355.0
$ws.Range("AJ51").Value = -65280.0
$ws.Range("G52").Value = 9190.0
$ws.Range("H52").Value = 0.22
$ws.Range("AI52").Value = -9190.0
$ws.Range("AJ52").Value = -74470.0
$ws.Range("G53").Value = 10109.0
$ws.Range("H53").Value = 0.22
$ws.Range("AI53").Value = -10109.0
$ws.Range("AJ53").Value = -84579.0
$ws.Range("G54").Value = 11120.0
$ws.Range("H54").Value = 0.22
$ws.Range("AI54").Value = -11120.0
$ws.Range("AJ54").Value = -95699.0
$ws.Range("G55").Value = 12232.0
$ws.Range("H55").Value = 0.22
$ws.Range("AI55").Value = -12232.0
$ws.Range("AJ55").Value = -107931.0
$ws.Range("G56").Value = 13455.0
$ws.Range("H56").Value = 0.22
$ws.Range("AI56").Value = -13455.0
$ws.Range("AJ56").Value = -121386.0
$ws.Range("G57").Value = 14801.0
$ws.Range("H57").Value = 0.22
$ws.Range("AI57").Value = -14801.0
$ws.Range("AJ57").Value = -136187.0
$ws.Range("G58").Value = 16281.0
$ws.Range("H58").Value = 0.22
$ws.Range("AI58").Value = -16281.0
$ws.Range("AJ58").Value = -152468.0
$ws.Range("G59").Value = 17909.0
$ws.Range("H59").Value = 0.22
$ws.Range("AI59").Value = -17909.0
$ws.Range("AJ59").Value = -170377.0
$ws.Range("G60").Value = 19699.0
$ws.Range("H60").Value = 0.22
$ws.Range("AI60").Value = -19699.0
$ws.Range("AJ60").Value = -190076.0
$ws.Range("G61").Value = 21670.0
$ws.Range("H61").Value = 0.22
$ws.Range("AI61").Value = -21670.0
$ws.Range("AJ61").Value = -211746.0
$ws.Range("G62").Value = 23836.0
$ws.Range("H62").Value = 0.22
$ws.Range("AI62").Value = -23836.0
$ws.Range("AJ62").Value = -235582.0
$ws.Range("G63").Value = 26220.0
$ws.Range("H63").Value = 0.22
$ws.Range("AI63").Value = -26220.0
$ws.Range("AJ63").Value = -261802.0
$ws.Range("G64").Value = 28842.0
$ws.Range("H64").Value = 0.22
$ws.Range("AI64").Value = -28842.0
$ws.Range("AJ64").Value = -290644.0
$ws.Range("G65").Value = 31726.0
$ws.Range("H65").Value = 0.22
$ws.Range("AI65").Value = -31726.0
$ws.Range("AJ65").Value = -322370.0
$ws.Range("G66").Value = 34899.0
$ws.Range("H66").Value = 0.22
$ws.Range("AI66").Value = -35348.44
$ws.Range("AJ66").Value = -357718.44
$ws.Range("G67").Value = 38389.0
$ws.Range("H67").Value = 0.22
$ws.Range("AI67").Value = -40583.38
$ws.Range("AJ67").Value = -398301.82
$ws.Range("G68").Value = 42228.0
$ws.Range("H68").Value = 0.22
$ws.Range("AI68").Value = -46341.81
$ws.Range("AJ68").Value = -444643.63
$ws.Range("G69").Value = 46450.0
$ws.Range("H69").Value = 0.22
$ws.Range("AI69").Value = -52675.2
$ws.Range("AJ69").Value = -497318.83
$ws.Range("G70").Value = 51095.0
$ws.Range("H70").Value = 0.22
$ws.Range("AI70").Value = -59642.72
$ws.Range("AJ70").Value = -556961.55
$ws.Range("G71").Value = 56205.0
$ws.Range("H71").Value = 0.22
$ws.Range("AI71").Value = -67307.49
$ws.Range("AJ71").Value = -624269.04
$ws.Range("G72").Value = 61826.0
$ws.Range("H72").Value = 0.22
$ws.Range("AI72").Value = -75738.73
$ws.Range("AJ72").Value = -700007.77
